$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text-format to all target cells first, then set values,
# so numeric-looking strings (e.g. "1.00", "0.539") are preserved exactly as text
$cellValues = @{
    'D2' = '50.871.32'
    'E2' = '  -1.46%  '
    'D3' = '2.903.79'
    'E3' = '  -1.43%  '
    'D4' = '1.00'
    'E4' = '  +0.02%  '
    'D5' = '367.80'
    'E5' = '  +4.88%  '
    'D6' = '102.27'
    'E6' = '  -3.76%  '
    'D7' = '0.539'
    'E7' = '  -2.73%  '
    'D8' = '1.00'
    'E8' = '  +0.04%  '
    'D9' = '0.580'
    'E9' = '  -3.95%  '
    'D10' = '36.71'
    'E10' = '  -2.90%  '
    'E11' = '  +0.71%  '
    'D12' = '0.0831'
    'E12' = '  -1.95%  '
    'D13' = '18.21'
    'E13' = '  -3.90%  '
    'D14' = '3.360.45'
    'E14' = '  -1.10%  '
    'E15' = '  -2.98%  '
    'D16' = '2.905.19'
    'E16' = '  -1.19%  '
    'D17' = '0.920'
    'E17' = '  -4.95%  '
    'D18' = '50.869.43'
    'E18' = '  -1.35%  '
    'D19' = '3.19'
    'E19' = '  -5.43%  '
    'D20' = '7.15'
    'E20' = '  -3.24%  '
    'D21' = '12.82'
    'E21' = '  -4.41%  '
    'D22' = '0.0₃0938'
    'E22' = '  -2.70%  '
    'D23' = '67.85'
    'E23' = '  -1.63%  '
    'D24' = '257.54'
    'E24' = '  -1.33%  '
    'E25' = '  -1.57%  '
    'D26' = '4.19'
    'E26' = '  -2.02%  '
    'E28' = '  -4.38%  '
    'D29' = '25.49'
    'E29' = '  -3.80%  '
    'D30' = '7.06'
    'E30' = '  -2.99%  '
    'E31' = '  -4.71%  '
    'D32' = '6.23'
    'E32' = '  +2.76%  '
    'D33' = '9.84'
    'E33' = '  -3.69%  '
    'E34' = '  -3.87%  '
    'D35' = '51.27'
    'E35' = '  +1.55%  '
    'D36' = '34.02'
    'E36' = '  -4.54%  '
    'E37' = '  +0.60%  '
    'D38' = '0.0419'
    'E38' = '  -2.58%  '
    'D39' = '2.96'
    'E39' = '  -5.68%  '
    'D40' = '16.96'
    'E40' = '  -4.08%  '
    'D41' = '2.59'
    'E41' = '  -2.41%  '
    'D42' = '1.83'
    'E42' = '  -5.68%  '
    'E43' = '  -2.99%  '
    'B44' = 'Monero'
    'C44' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D44' = '119.12'
    'E44' = '  -4.03%  '
    'B45' = 'EnergySwap'
    'C45' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D45' = '21.86'
    'E45' = '  -1.65%  '
    'E46' = '  -1.92%  '
    'D47' = '2.014.81'
    'E47' = '  -4.12%  '
    'E48' = '  -0.03%  '
    'D49' = '3.12'
    'E49' = '  -5.80%  '
    'D50' = '3.188.43'
    'D51' = '0.235'
    'E51' = '  -0.24%  '
}

foreach ($cellRef in $cellValues.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $cellValues[$cellRef]
}